$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to carry two header rows: row 1 with unit labels in a few
# columns, and row 2 with "Hiver/Eté/Année" season labels. Remove row 2
# entirely (all data rows shift up by one) and turn row 1 into a single,
# complete header row that also introduces the new idx/idx2/Name/Date
# Start/Date End columns plus renamed power & energy columns.
$ws.Rows(2).Delete()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 keep the plain default style, F1:K1 take on the (font-only) header
# style used across the rest of the header row.
$ws.Range("A1:E1").Font.Size = 10
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").NumberFormat = "General"

$ws.Range("A2:K2").Select()
